$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.883.01"
$ws.Range("E2").Value = "  +0.82%  "

$ws.Range("D3").Value = "3.110.59"
$ws.Range("E3").Value = "  +4.29%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "390.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.24%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.547"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.34%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.591"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.17%  "

$ws.Range("E10").Value = "  +1.44%  "

$ws.Range("E11").Value = "  +0.47%  "

$ws.Range("E12").Value = "  +1.12%  "

$ws.Range("D13").Value = "3.579.42"
$ws.Range("E13").Value = "  +3.79%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.97%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.02%  "

$ws.Range("D16").Value = "3.089.81"
$ws.Range("E16").Value = "  +3.86%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.985"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.30%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.14%  "

$ws.Range("D19").Value = "51.949.11"
$ws.Range("E19").Value = "  +0.96%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.25%  "

$ws.Range("D22").Value = "0.0₃0973"
$ws.Range("E22").Value = "  +1.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.26%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "270.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.22%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.86%  "

$ws.Range("E28").Value = "  +2.25%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.09%  "

$ws.Range("E30").Value = "  +0.15%  "

$ws.Range("E31").Value = "  -0.63%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.30%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.71"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.77%  "

$ws.Range("E34").Value = "  +0.54%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.22%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.41"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.88%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.298"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.08%  "

$ws.Range("E40").Value = "  +2.69%  "

$ws.Range("E41").Value = "  +1.17%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.59"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.81%  "

$ws.Range("E43").Value = "  -0.18%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "126.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.43%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.91%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.48%  "

$ws.Range("E47").Value = "  +4.31%  "

$ws.Range("E48").Value = "  +3.07%  "

$ws.Range("D49").Value = "2.046.57"
$ws.Range("E49").Value = "  +0.93%  "

$ws.Range("D50").Value = "3.401.00"
$ws.Range("E50").Value = "  +3.67%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.209"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.76%  "

# Rows 35/36: OKB and VeChain swap positions (re-sorted by price), with
# refreshed price/volume figures.
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.36"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.10%  "

$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0452"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.42%  "

